$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -0.5664960695472127
$ws.Range("B3").Value = -0.4361400413227467
$ws.Range("B4").Value = -0.5371685250986227
$ws.Range("B5").Value = -0.4931828166462521
$ws.Range("B6").Value = -0.5571679912696158
$ws.Range("B7").Value = -0.5466931779871508
$ws.Range("B8").Value = -0.5712379502829247
$ws.Range("B9").Value = -0.7974733395810489
$ws.Range("B10").Value = -0.5425470077601648
$ws.Range("B11").Value = -0.4026775245560015
$ws.Range("B12").Value = -0.3404220649242103
$ws.Range("B13").Value = -0.546406507522678
$ws.Range("B14").Value = 0.4144728952908346
$ws.Range("B15").Value = 0.02358490778145368
$ws.Range("B16").Value = -0.05931516769736236
$ws.Range("B17").Value = 0.2266109970553055
$ws.Range("B18").Value = -0.2077798954074883
$ws.Range("B19").Value = -0.3069793719305974
$ws.Range("B20").Value = 0.01155715253908448
$ws.Range("B21").Value = 0.284731672537188
$ws.Range("B22").Value = -0.005243094914738989
$ws.Range("B23").Value = 0.06392999676223968
$ws.Range("B24").Value = 0.09990703194902172
$ws.Range("B25").Value = -0.2512958655336935
$ws.Range("B26").Value = -0.1658123980933962
$ws.Range("B27").Value = -0.03741110982631637
$ws.Range("B28").Value = 0.2179162127468129
$ws.Range("B29").Value = -0.1275363987006121
$ws.Range("B30").Value = -0.09421348167506201
$ws.Range("B31").Value = 0.118423199919853
$ws.Range("B32").Value = -0.1920792050671807
$ws.Range("B33").Value = -0.1776041784525154
$ws.Range("B34").Value = -0.0727936468978014
$ws.Range("B35").Value = -0.01510879135064917
$ws.Range("B36").Value = -0.1994933507736878
$ws.Range("B37").Value = 0.23605657102707
$ws.Range("B38").Value = -0.3022795516967504
$ws.Range("B39").Value = -0.6571408764479121
$ws.Range("B40").Value = -0.2823974595631665
$ws.Range("B41").Value = -0.1753352089740312
$ws.Range("B42").Value = -0.3147999646388679
$ws.Range("B43").Value = -0.1847239898599183
$ws.Range("B44").Value = -0.2340311456644844
$ws.Range("B45").Value = -0.2412528006065553
$ws.Range("B46").Value = -0.2751932964850113
$ws.Range("B47").Value = -0.6284685788168278
$ws.Range("B48").Value = -0.374323726930951
$ws.Range("B49").Value = -0.2519311507540934
$ws.Range("B50").Value = 0.07433946178030867
$ws.Range("B51").Value = -0.1366867419062487
$ws.Range("B52").Value = -0.01172748182599049
$ws.Range("B53").Value = -0.002100530557000856
$ws.Range("B54").Value = -0.0009408250518614793
$ws.Range("B55").Value = 0.1976166453420277
$ws.Range("B56").Value = -0.009501117240594794
$ws.Range("B57").Value = -0.01564762683657516
$ws.Range("B58").Value = -0.0414926571274016
$ws.Range("B59").Value = 0.006632197058229857
$ws.Range("B60").Value = 0.01688795277747085
$ws.Range("B61").Value = 0.007894357705547264
$ws.Range("B62").Value = 0.1296776238540422
$ws.Range("B63").Value = 0.2130853869571952
$ws.Range("B64").Value = -0.203853777179843
$ws.Range("B65").Value = -0.008382073170275468
$ws.Range("B66").Value = -0.0143028176851671
$ws.Range("B67").Value = 0.03231732722739699
$ws.Range("B68").Value = 0.01093885602750555
$ws.Range("B69").Value = 0.02703306490989708
$ws.Range("B70").Value = 0.2455166095388436
$ws.Range("B71").Value = -0.07134232224975101
$ws.Range("B72").Value = 0.4310516329957331
$ws.Range("B73").Value = 0.2263222698669136
$ws.Range("B74").Value = 0.2093110868181786
$ws.Range("B75").Value = 0.4367228708906591
$ws.Range("B76").Value = 0.2066625670401257
$ws.Range("B77").Value = -0.2089765009323316
$ws.Range("B78").Value = -0.08127697723661732
$ws.Range("B79").Value = 0.4357946272645412
$ws.Range("B80").Value = 0.02589980063020375
$ws.Range("B81").Value = 0.1986133481052239
$ws.Range("B82").Value = 0.1303470659471042
$ws.Range("B83").Value = 0.01094085719139458
$ws.Range("B84").Value = 0.173748220711727
$ws.Range("B85").Value = 0.2504041201479615
$ws.Range("B86").Value = 0.1202593546188853
$ws.Range("B87").Value = 0.2675351494061741
$ws.Range("B88").Value = 0.3570341646186974
$ws.Range("B89").Value = 0.04941819241653365
$ws.Range("B90").Value = 0.3670699540512709
$ws.Range("B91").Value = 0.272985159173188
$ws.Range("B92").Value = 0.2312854247614052
$ws.Range("B93").Value = 0.03638017565047692
$ws.Range("B94").Value = 0.1161180873737637
$ws.Range("B95").Value = 0.05748720845702045
$ws.Range("B96").Value = -0.0215153756233217
$ws.Range("B97").Value = 0.1617353827136256
$ws.Range("B98").Value = 0.5533821121420188
$ws.Range("B99").Value = 0.484205464036842
$ws.Range("B100").Value = 0.4236925397683906
$ws.Range("B101").Value = 0.7454669809451054
$ws.Range("B102").Value = 0.3310542140074799
$ws.Range("B103").Value = 0.2631345108208606
$ws.Range("B104").Value = 0.3185023476819153
$ws.Range("B105").Value = 0.0979728812981268
$ws.Range("B106").Value = 0.449268492969307
$ws.Range("B107").Value = 0.6822795084745374
$ws.Range("B108").Value = 0.502950008763768
$ws.Range("B109").Value = 0.2975128394246513
$ws.Range("B110").Value = 0.2694575730101699
$ws.Range("B111").Value = -0.07754066984257242
$ws.Range("B112").Value = -0.0477600149143873
$ws.Range("B113").Value = 0.1755379928819636
$ws.Range("B114").Value = -0.1004729398308437
$ws.Range("B115").Value = 0.162615240589367
$ws.Range("B116").Value = 0.1038539838732898
$ws.Range("B117").Value = 0.1141277493459276
$ws.Range("B118").Value = 0.1057233991064729
$ws.Range("B119").Value = -0.07370327716450033
$ws.Range("B120").Value = -0.06244402007442865
$ws.Range("B121").Value = 0.2121881080823664
$ws.Range("B122").Value = -0.3764873681567382
$ws.Range("B123").Value = -0.373037751830464
$ws.Range("B124").Value = -0.3978231460803199
$ws.Range("B125").Value = -0.4605336321505669
$ws.Range("B126").Value = -0.4380338817586157
$ws.Range("B127").Value = -0.4731803917823798
$ws.Range("B128").Value = -0.3206905203920072
$ws.Range("B129").Value = -0.4985917344926163
$ws.Range("B130").Value = -0.5868984659402303
$ws.Range("B131").Value = -0.3179646442813064
$ws.Range("B132").Value = -0.6394334570414045
$ws.Range("B133").Value = -0.3912465695228383
$ws.Range("B134").Value = 0.4905996676788515
$ws.Range("B135").Value = 0.4025874100312921
$ws.Range("B136").Value = 0.3305835078702661
$ws.Range("B137").Value = 0.3988252963059269
$ws.Range("B138").Value = 0.146418781965773
$ws.Range("B139").Value = 0.4076768932778474
$ws.Range("B140").Value = 0.132698968441055
$ws.Range("B141").Value = 0.2729028605682403
$ws.Range("B142").Value = 0.4783126136697041
$ws.Range("B143").Value = 0.59175860304156
$ws.Range("B144").Value = 0.5909453997744539
$ws.Range("B145").Value = 0.4844257811259522
